$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two match rows (23 and 24) that were previously in the wrong
# --- order (Murcia-Sanluqueno vs Linares-Real Madrid B). Keep the shared
# --- columns (Indice, pais, torneio, temporada, data_partida) as-is and
# --- only swap the per-match columns F..V.
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $swapCols) {
    $cellTop = $col + "23"
    $cellBottom = $col + "24"
    $valTop = $ws.Range($cellTop).Value()
    $valBottom = $ws.Range($cellBottom).Value()
    $ws.Range($cellTop).Value = $valBottom
    $ws.Range($cellBottom).Value = $valTop
}

# --- Append 4 new match rows (105-108) at the bottom of the sheet, copying
# --- the formatting of the last existing row (104) so the new rows line up
# --- with the sheet's existing cell styles (bold/centered/bordered Indice
# --- column, date-formatted data_partida column).
$ws.Range("A104:V104").Copy()
$ws.Range("A105:V108").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @(104, "spain", "primera-rfef-group-2", "2023-2024", 45235.66666666666, "Ceuta", 1, "Murcia", 0, 2.32, "02/11/2023 08:13", 2.08, "05/11/2023 15:36", 2.92, "02/11/2023 08:13", 3.04, "05/11/2023 15:36", 3.05, "02/11/2023 08:13", 4.05, "05/11/2023 15:36", "https://www.betexplorer.com/football/spain/primera-rfef-group-2/ad-ceuta-murcia/GKOhoj8B/"),
    @(105, "spain", "primera-rfef-group-2", "2023-2024", 45235.75, "Castellon", 2, "Alcoyano", 0, 1.53, "02/11/2023 08:13", 1.41, "05/11/2023 17:39", 3.73, "02/11/2023 08:13", 4.38, "05/11/2023 17:39", 5.48, "02/11/2023 08:13", 8.41, "05/11/2023 17:39", "https://www.betexplorer.com/football/spain/primera-rfef-group-2/castellon-alcoyano/fBNdpANH/"),
    @(106, "spain", "primera-rfef-group-2", "2023-2024", 45235.83333333334, "Malaga", 1, "Cordoba", 1, 1.98, "02/11/2023 08:13", 1.83, "05/11/2023 18:04", 3.47, "02/11/2023 08:13", 3.53, "05/11/2023 18:04", 3.36, "02/11/2023 08:13", 4.36, "05/11/2023 18:04", "https://www.betexplorer.com/football/spain/primera-rfef-group-2/malaga-cordoba/Ct7A3XNu/"),
    @(107, "spain", "primera-rfef-group-2", "2023-2024", 45235.83333333334, "Merida AD", 1, "UD Ibiza", 2, 3.88, "02/11/2023 08:13", 4.12, "05/11/2023 09:42", 3.19, "02/11/2023 08:13", 3.34, "05/11/2023 18:02", 1.88, "02/11/2023 08:13", 1.92, "05/11/2023 09:42", "https://www.betexplorer.com/football/spain/primera-rfef-group-2/merida-ad-ud-ibiza/K6AI1Bhh/")
)

$startRow = 105
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
